$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "document.contract" row (row 7), cloning row 6's formatting ---
$ws.Range("A6:I6").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A7").Value = "6cdd34ca-1e4d-4e63-ba78-b0ead49a8fd2"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "zero.document.default"
$ws.Range("D7").Value = "合同管理"
$ws.Range("E7").Value = "document.contract"
$ws.Range("F7").Value = $true
$ws.Range("H7").Value = 1015

# --- Update the view's selection / scroll position ---
[void]$ws.Range("H8").Select()
